# "Full automation complete #3" - tidy up the header row on every per-stock
# sheet: drop the verbose "XXX Price ($)" / "Time (hh:mm:ss)" / "Current
# Price ($)" wording in favour of shorter "Time" / "Open ($)" / "High ($)" /
# "Low ($)" / "Close ($)" / "Price ($)" headers, reorder the columns (Time,
# Open, High, Low, Close, Price, Change($), Change(%)), re-select the last
# worked cell on each tab, and bring the Monster tab back to the front.

$wb = $excel.ActiveWorkbook

$wsMonster   = $wb.Worksheets.Item("Monster")
$wsSamsung   = $wb.Worksheets.Item("Samsung")
$wsEA        = $wb.Worksheets.Item("EA")
$wsKonami    = $wb.Worksheets.Item("Konami")
$wsNvidia    = $wb.Worksheets.Item("Nvidia")
$wsMicrosoft = $wb.Worksheets.Item("Microsoft")

# --- Header row tidy-up -----------------------------------------------
# Column A + F (and G/H, already fine) first across every sheet so the new
# shared-string table fills in "Time" / "Price ($)" before the "Time "
# (trailing-space) variant, matching how the strings were actually typed in
# while going sheet-by-sheet during cleanup.

$wsKonami.Range("A3").Value = "Time"
$wsKonami.Range("F3").Value = "Price ($)"
$wsKonami.Range("G3").Value = "Change ($)"
$wsKonami.Range("H3").Value = "Change (%)"

$wsMicrosoft.Range("A3").Value = "Time"
$wsMicrosoft.Range("F3").Value = "Price ($)"
$wsMicrosoft.Range("G3").Value = "Change ($)"
$wsMicrosoft.Range("H3").Value = "Change (%)"

$wsMonster.Range("A3").Value = "Time "
$wsMonster.Range("F3").Value = "Price ($)"
$wsMonster.Range("G3").Value = "Change ($)"
$wsMonster.Range("H3").Value = "Change (%)"

$wsEA.Range("A3").Value = "Time "
$wsEA.Range("F3").Value = "Price ($)"
$wsEA.Range("G3").Value = "Change ($)"
$wsEA.Range("H3").Value = "Change (%)"

$wsNvidia.Range("A3").Value = "Time "
$wsNvidia.Range("F3").Value = "Price ($)"
$wsNvidia.Range("G3").Value = "Change ($)"
$wsNvidia.Range("H3").Value = "Change (%)"

$wsKonami.Range("B3").Value = "Open ($)"
$wsKonami.Range("C3").Value = "High ($)"
$wsKonami.Range("D3").Value = "Low ($)"
$wsKonami.Range("E3").Value = "Close ($)"

$wsMicrosoft.Range("B3").Value = "Open ($)"
$wsMicrosoft.Range("C3").Value = "High ($)"
$wsMicrosoft.Range("D3").Value = "Low ($)"
$wsMicrosoft.Range("E3").Value = "Close ($)"

$wsMonster.Range("B3").Value = "Open ($)"
$wsMonster.Range("C3").Value = "High ($)"
$wsMonster.Range("D3").Value = "Low ($)"
$wsMonster.Range("E3").Value = "Close ($)"

$wsEA.Range("B3").Value = "Open ($)"
$wsEA.Range("C3").Value = "High ($)"
$wsEA.Range("D3").Value = "Low ($)"
$wsEA.Range("E3").Value = "Close ($)"

$wsNvidia.Range("B3").Value = "Open ($)"
$wsNvidia.Range("C3").Value = "High ($)"
$wsNvidia.Range("D3").Value = "Low ($)"
$wsNvidia.Range("E3").Value = "Close ($)"

# --- Microsoft header row: clear the leftover wrap-text flag ----------
# Leaves an (empty) alignment record behind, same as toggling WrapText off
# in the UI after it had been fiddled with.
$wsMicrosoft.Range("A3:H3").WrapText = $false

# --- Page setup tidy-up (EA + Microsoft) -------------------------------
$wsEA.PageSetup.PaperSize = 9
$wsEA.PageSetup.Orientation = 1

$wsMicrosoft.PageSetup.PaperSize = 9
$wsMicrosoft.PageSetup.Orientation = 1

# --- Re-apply the last-worked selection on every tab -------------------
[void]$wsSamsung.Activate()
[void]$wsSamsung.Range("L22").Select()

[void]$wsEA.Activate()
[void]$wsEA.Range("G16").Select()

[void]$wsKonami.Activate()
[void]$wsKonami.Range("G10").Select()

[void]$wsNvidia.Activate()
[void]$wsNvidia.Range("F11").Select()

[void]$wsMicrosoft.Activate()
[void]$wsMicrosoft.Range("F5").Select()

# Monster ends up on top, as the front/active tab.
[void]$wsMonster.Activate()
[void]$wsMonster.Range("F5").Select()
